$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at 271, shifting the existing rows 271-336 down to 272-337.
$ws.Rows.Item(271).Insert()

# Populate the newly inserted row 271 with its data (matches target diff exactly).
$ws.Cells.Item(271, 1).Value = 9
$ws.Cells.Item(271, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(271, 3).Value = "Metropolitana"
$ws.Cells.Item(271, 4).Value = 44964
$ws.Cells.Item(271, 5).Value = 13
$ws.Cells.Item(271, 6).Value = 100112001
$ws.Cells.Item(271, 7).Value = "Berenjena"
$ws.Cells.Item(271, 8).Value = "Sin especificar"
$ws.Cells.Item(271, 9).Value = "Primera"
$ws.Cells.Item(271, 10).Value = 52
$ws.Cells.Item(271, 11).Value = 10000
$ws.Cells.Item(271, 12).Value = 11000
$ws.Cells.Item(271, 13).Value = 10500
$ws.Cells.Item(271, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(271, 15).Value = "Región Metropolitana"
$ws.Cells.Item(271, 16).Value = 210
$ws.Cells.Item(271, 17).Value = 50
$ws.Cells.Item(271, 18).Value = "Hortaliza"
